$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Checking for Anagrams - Hash"
$ws.Range("H16").Value = "CheckAnagrams"
$ws.Range("C16").Value = "Hash related problems"
$null = $ws.Range("C17").Select()
